# ---------------------------------------------------------------------------
# Update countries & provincias Spain
#
# Refreshes the COVID-19 stats snapshot on sheet "Pais":
#   - updates the "Datos actualizados..." timestamp caption (A1)
#   - updates Casos totales/Nuevos casos/Casos activos/Recuperados/
#     Casos criticos/Muertes hoy/Muertes for the countries whose figures
#     changed in this refresh
#   - because several countries changed rank (table is sorted by Casos
#     totales, column B, descending), a few rows end up showing a new
#     country name even though only neighbouring rows totals moved
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Header caption timestamp
$ws.Range("A1").Value = 'Datos actualizados a 12 de Abril de 2020 a las 06:52'

# Row 17: Brasil
$ws.Range("B17").Value = 20964
$ws.Range("C17").Value = 2
$ws.Range("E17").Value = 19650
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 1141

# Row 32: Australia
$ws.Range("E32").Value = 2979
$ws.Range("G32").Value = 3
$ws.Range("H32").Value = 59

# Row 50: Tailandia
$ws.Range("B50").Value = 2551
$ws.Range("C50").Value = 33
$ws.Range("D50").Value = 1218
$ws.Range("E50").Value = 1295
$ws.Range("G50").Value = 3
$ws.Range("H50").Value = 38

# Row 102: Kirguistan
$ws.Range("A102").Value = 'Kirguistan'
$ws.Range("B102").Value = 377
$ws.Range("C102").Value = 38
$ws.Range("D102").Value = 44
$ws.Range("E102").Value = 328
$ws.Range("F102").Value = 5
$ws.Range("H102").Value = 5

# Row 103: Malta
$ws.Range("A103").Value = 'Malta'
$ws.Range("B103").Value = 370
$ws.Range("D103").Value = 16
$ws.Range("E103").Value = 351
$ws.Range("F103").Value = 4
$ws.Range("H103").Value = 3

# Row 104: San Marino
$ws.Range("A104").Value = 'San Marino'
$ws.Range("B104").Value = 356
$ws.Range("D104").Value = 53
$ws.Range("E104").Value = 268
$ws.Range("F104").Value = 14
$ws.Range("H104").Value = 35

# Row 165: Antigua y Barbuda
$ws.Range("A165").Value = 'Antigua y Barbuda'
$ws.Range("D165").Value = 0
$ws.Range("F165").Value = 1
$ws.Range("H165").Value = 2

# Row 166: Somalia
$ws.Range("A166").Value = 'Somalia'
$ws.Range("D166").Value = 2
$ws.Range("E166").Value = 18
$ws.Range("F166").Value = 2
$ws.Range("H166").Value = 1

# Row 187: Nepal
$ws.Range("A187").Value = 'Nepal'
$ws.Range("C187").Value = 3
$ws.Range("D187").Value = 1
$ws.Range("E187").Value = 11
$ws.Range("F187").Value = 0
$ws.Range("H187").Value = 0

# Row 188: Malaui
$ws.Range("A188").Value = 'Malaui'
$ws.Range("D188").Value = 0
$ws.Range("E188").Value = 10
$ws.Range("F188").Value = 1
$ws.Range("H188").Value = 2

# Row 189: Suazilandia
$ws.Range("A189").Value = 'Suazilandia'
$ws.Range("B189").Value = 12
$ws.Range("D189").Value = 7
$ws.Range("E189").Value = 5

# Row 190: Seychelles
$ws.Range("A190").Value = 'Seychelles'
$ws.Range("D190").Value = 0
$ws.Range("E190").Value = 11

# Row 191: Republica del Chad
$ws.Range("A191").Value = 'Republica del Chad'
$ws.Range("D191").Value = 2
$ws.Range("E191").Value = 9

# Row 192: Groenlandia
$ws.Range("A192").Value = 'Groenlandia'
$ws.Range("B192").Value = 11
$ws.Range("D192").Value = 11
$ws.Range("E192").Value = 0

# Row 193: Sierra Leona
$ws.Range("A193").Value = 'Sierra Leona'
$ws.Range("D193").Value = 0
$ws.Range("E193").Value = 10

# Row 194: Surinam
$ws.Range("A194").Value = 'Surinam'
$ws.Range("B194").Value = 10
$ws.Range("D194").Value = 4
$ws.Range("E194").Value = 5
$ws.Range("H194").Value = 0

# Row 195: Islas Turcas y Caicos
$ws.Range("A195").Value = 'Islas Turcas y Caicos'
$ws.Range("D195").Value = 0
$ws.Range("H195").Value = 1

# Row 196: Nicaragua
$ws.Range("A196").Value = 'Nicaragua'
